$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Strip mojibake trademark/accent artifacts from menu item names, and relabel the
# classification headers from letter codes to descriptive labels. Writes are sequenced
# so each distinct cleaned-up string is introduced (first use) in a specific order.
$ws.Range("B3").Value = 'Double Quarter Pounder with Cheese'
$ws.Range("B4").Value = 'Quarter Pounder with Cheese'
$ws.Range("B5").Value = 'Big Mac'
$ws.Range("B7").Value = 'Chicken McCrispy (2pc)'
$ws.Range("B10").Value = 'Double McSpicy'
$ws.Range("B11").Value = 'McSpicy'
$ws.Range("B12").Value = 'Double Filet-O-Fish'
$ws.Range("B13").Value = 'Filet-O-Fish'
$ws.Range("B14").Value = 'Chicken McNuggets (9pc)'
$ws.Range("B15").Value = 'Chicken McNuggets (6pc)'
$ws.Range("B16").Value = 'McWings (4pc)'
$ws.Range("B17").Value = 'McChicken'
$ws.Range("B18").Value = 'Grilled Chicken McWrap'
$ws.Range("B48").Value = 'Chicken McCrispy (6pc)'
$ws.Range("B65").Value = 'Chicken McNuggets (20pc)'
$ws.Range("B69").Value = 'McWings (2pc)'
$ws.Range("B71").Value = 'OREO McFlurry'
$ws.Range("B72").Value = 'Mudpie McFlurry'
$ws.Range("B73").Value = 'Strawberry Shortcake McFlurry'
$ws.Range("B76").Value = 'ChocoCone'
$ws.Range("B82").Value = 'Iced MILO (Small)'
$ws.Range("B92").Value = 'Hot MILO'
$ws.Range("B93").Value = 'Coca-Cola Zero Sugar (Small)'
$ws.Range("B94").Value = 'Coca-Cola Original Taste Less Sugar (Small)'
$ws.Range("B95").Value = 'Sprite (Small)'
$ws.Range("B102").Value = 'Big Breakfast'
$ws.Range("B107").Value = 'Sausage McMuffin'
$ws.Range("B108").Value = 'Sausage McMuffin with Egg'
$ws.Range("B111").Value = 'Egg McMuffin'
$ws.Range("B80").Value = 'Caramel Frappe (Small)'
$ws.Range("B81").Value = 'Mocha Frappe (Small)'
$ws.Range("B85").Value = 'McCafe Cappuccino'
$ws.Range("B86").Value = 'McCafe Latte'
$ws.Range("B87").Value = 'McCafe Iced Latte'
$ws.Range("B88").Value = 'McCafe Premium Roast Coffee'
$ws.Range("D1").Value = 'Low'
$ws.Range("E1").Value = 'Medium'
$ws.Range("F1").Value = 'High'

# --- Remaining repeated occurrences of the same cleaned-up menu item names ---
$ws.Range("B8").Value = 'Chicken McCrispy (2pc)'
$ws.Range("B23").Value = 'Double Quarter Pounder with Cheese'
$ws.Range("B24").Value = 'Quarter Pounder with Cheese'
$ws.Range("B25").Value = 'Big Mac'
$ws.Range("B27").Value = 'Chicken McCrispy (2pc)'
$ws.Range("B28").Value = 'Chicken McCrispy (2pc)'
$ws.Range("B30").Value = 'Double McSpicy'
$ws.Range("B31").Value = 'McSpicy'
$ws.Range("B32").Value = 'Double Filet-O-Fish'
$ws.Range("B33").Value = 'Filet-O-Fish'
$ws.Range("B34").Value = 'Chicken McNuggets (9pc)'
$ws.Range("B35").Value = 'Chicken McNuggets (6pc)'
$ws.Range("B36").Value = 'McWings (4pc)'
$ws.Range("B37").Value = 'McChicken'
$ws.Range("B38").Value = 'Grilled Chicken McWrap'
$ws.Range("B43").Value = 'Double Quarter Pounder with Cheese'
$ws.Range("B44").Value = 'Quarter Pounder with Cheese'
$ws.Range("B45").Value = 'Big Mac'
$ws.Range("B47").Value = 'Chicken McCrispy (2pc)'
$ws.Range("B50").Value = 'Double McSpicy'
$ws.Range("B51").Value = 'McSpicy'
$ws.Range("B52").Value = 'Double Filet-O-Fish'
$ws.Range("B53").Value = 'Filet-O-Fish'
$ws.Range("B54").Value = 'Chicken McNuggets (9pc)'
$ws.Range("B55").Value = 'Chicken McNuggets (6pc)'
$ws.Range("B57").Value = 'McChicken'
$ws.Range("B58").Value = 'Grilled Chicken McWrap'
$ws.Range("B64").Value = 'McWings (4pc)'
$ws.Range("B112").Value = 'Double Filet-O-Fish'
$ws.Range("B113").Value = 'Filet-O-Fish'
$ws.Range("B119").Value = 'Big Breakfast'
$ws.Range("B124").Value = 'Sausage McMuffin with Egg'
$ws.Range("B125").Value = 'Sausage McMuffin'
$ws.Range("B128").Value = 'Egg McMuffin'

# --- Move the active selection ---
$ws.Range("D2").Select()
